$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3088.7778
$ws.Range("I18").Value = 2960.2
$ws.Range("K18").Value = 2960.2
$ws.Range("M18").Value = -2676.2

$ws.Range("H19").Value = 4717.8667
$ws.Range("I19").Value = 2350.4
$ws.Range("K19").Value = 2350.4
$ws.Range("M19").Value = -2175.4

$ws.Range("H53").Value = 482.2439
$ws.Range("I53").Value = 289.15152
$ws.Range("J53").Value = 1278.75
$ws.Range("K53").Value = 289.15152
$ws.Range("L53").Value = 1278.75
$ws.Range("M53").Value = 347.84848
$ws.Range("N53").Value = -2552.75

$ws.Range("H101").Value = 2410.8333
$ws.Range("J101").Value = 4500
$ws.Range("L101").Value = 13500
$ws.Range("N101").Value = -16744

$ws.Range("H138").Value = 3652.25
$ws.Range("I138").Value = 3214.4546
$ws.Range("J138").Value = 3830.611
$ws.Range("K138").Value = 9643.363799999999
$ws.Range("L138").Value = 11491.833
$ws.Range("M138").Value = -4503.363799999999
$ws.Range("N138").Value = -21771.833

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 17.5
$ws.Range("I5").Value = 17.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 17.5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 94.5
$ws.Range("N5").ClearContents()

$ws.Range("H32").Value = 23033.416
$ws.Range("I32").Value = 23107.613
$ws.Range("K32").Value = 23107.613
$ws.Range("M32").Value = -22820.613

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 17.5
$ws.Range("I4").Value = 17.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 17.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 97.5
$ws.Range("N4").ClearContents()

$ws.Range("H20").Value = 43211.28
$ws.Range("I20").Value = 52613.7
$ws.Range("K20").Value = 52613.7
$ws.Range("M20").Value = -52366.7

$ws.Range("H86").Value = 89404.35000000001
$ws.Range("I86").Value = 1583.4445
$ws.Range("K86").Value = 1583.4445
$ws.Range("M86").Value = -460.4445000000001

$ws.Range("H89").Value = 89404.35000000001
$ws.Range("I89").Value = 1583.4445
$ws.Range("K89").Value = 7917.2225
$ws.Range("M89").Value = -2301.2225

$ws.Range("H99").Value = 26176.875
$ws.Range("I99").Value = 27255.334
$ws.Range("K99").Value = 27255.334
$ws.Range("M99").Value = -25757.334

$ws.Range("H107").Value = 1480.6875
$ws.Range("I107").Value = 1311
$ws.Range("K107").Value = 1311
$ws.Range("M107").Value = 609

$ws.Range("H134").Value = 4976.3794
$ws.Range("I134").Value = 3724.5417
$ws.Range("J134").Value = 10985.2
$ws.Range("K134").Value = 11173.6251
$ws.Range("L134").Value = 32955.60000000001
$ws.Range("M134").Value = -8638.625100000001
$ws.Range("N134").Value = -38025.60000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1061.2222
$ws.Range("I22").Value = 1160.2
$ws.Range("J22").Value = 937.5
$ws.Range("K22").Value = 1160.2
$ws.Range("L22").Value = 937.5
$ws.Range("M22").Value = -810.2
$ws.Range("N22").Value = -1637.5

$ws.Range("H132").Value = 24947.209
$ws.Range("I132").Value = 1166.2
$ws.Range("K132").Value = 3498.6
$ws.Range("M132").Value = -968.6000000000004

$ws.Range("H134").Value = 3847.7273
$ws.Range("J134").Value = 5296.5713
$ws.Range("L134").Value = 15889.7139
$ws.Range("N134").Value = -20959.7139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1175
$ws.Range("J22").Value = 1500
$ws.Range("L22").Value = 4500
$ws.Range("N22").Value = -4838

$ws.Range("H27").Value = 1175
$ws.Range("J27").Value = 1500
$ws.Range("L27").Value = 4500
$ws.Range("N27").Value = -4704

$ws.Range("H39").Value = 2189.4
$ws.Range("J39").Value = 2599.25
$ws.Range("L39").Value = 7797.75
$ws.Range("N39").Value = -8385.75

$ws.Range("H55").Value = 3686.1667
$ws.Range("I55").Value = 529.75
$ws.Range("J55").Value = 9999
$ws.Range("K55").Value = 1589.25
$ws.Range("L55").Value = 29997
$ws.Range("M55").Value = -1412.25
$ws.Range("N55").Value = -30351

$ws.Range("H92").Value = 797.93335
$ws.Range("I92").Value = 1073.5555
$ws.Range("J92").Value = 384.5
$ws.Range("K92").Value = 3220.6665
$ws.Range("L92").Value = 1153.5
$ws.Range("M92").Value = -1972.6665
$ws.Range("N92").Value = -3649.5

$ws.Range("H122").Value = 55560372
$ws.Range("J122").Value = 952.6667
$ws.Range("L122").Value = 8574.0003
$ws.Range("N122").Value = -13474.0003

$ws.Range("H130").Value = 20000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 20000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 60000
$ws.Range("M130").ClearContents()
$ws.Range("N130").Value = -70040

$ws.Range("H131").Value = 11500406
$ws.Range("J131").Value = 8740.421
$ws.Range("L131").Value = 26221.263
$ws.Range("N131").Value = -36301.263

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5466.273
$ws.Range("I80").Value = 3953.2
$ws.Range("J80").Value = 6727.1665
$ws.Range("K80").Value = 3953.2
$ws.Range("L80").Value = 6727.1665
$ws.Range("M80").Value = -2955.2
$ws.Range("N80").Value = -8723.166499999999

$ws.Range("H83").Value = 5466.273
$ws.Range("I83").Value = 3953.2
$ws.Range("J83").Value = 6727.1665
$ws.Range("K83").Value = 19766
$ws.Range("L83").Value = 33635.8325
$ws.Range("M83").Value = -14774
$ws.Range("N83").Value = -43619.8325

$ws.Range("H111").Value = 29999
$ws.Range("J111").Value = 29999
$ws.Range("L111").Value = 29999
$ws.Range("N111").Value = -36133

$ws.Range("H113").Value = 4115.3335
$ws.Range("I113").Value = 3888.4
$ws.Range("J113").Value = 5250
$ws.Range("K113").Value = 3888.4
$ws.Range("L113").Value = 5250
$ws.Range("M113").Value = -1718.4
$ws.Range("N113").Value = -9590

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2249.3333
$ws.Range("I7").Value = 1999.2
$ws.Range("J7").Value = 3500
$ws.Range("K7").Value = 1999.2
$ws.Range("L7").Value = 3500
$ws.Range("M7").Value = -1887.2
$ws.Range("N7").Value = -3724

$ws.Range("H16").Value = 1408.6538
$ws.Range("I16").Value = 1021.5
$ws.Range("K16").Value = 1021.5
$ws.Range("M16").Value = -851.5

$ws.Range("H22").Value = 3615.147
$ws.Range("I22").Value = 1999.8462
$ws.Range("J22").Value = 4615.095
$ws.Range("K22").Value = 1999.8462
$ws.Range("L22").Value = 4615.095
$ws.Range("M22").Value = -1704.8462
$ws.Range("N22").Value = -5205.095

$ws.Range("H27").Value = 3615.147
$ws.Range("I27").Value = 1999.8462
$ws.Range("J27").Value = 4615.095
$ws.Range("K27").Value = 1999.8462
$ws.Range("L27").Value = 4615.095
$ws.Range("M27").Value = -1892.8462
$ws.Range("N27").Value = -4829.095

$ws.Range("H61").Value = 4911.05
$ws.Range("I61").Value = 4959.0527
$ws.Range("K61").Value = 4959.0527
$ws.Range("M61").Value = -4757.0527

$ws.Range("H82").Value = 2424.875
$ws.Range("I82").Value = 1999.8334
$ws.Range("K82").Value = 1999.8334
$ws.Range("M82").Value = -1638.8334

$ws.Range("H85").Value = 2424.875
$ws.Range("I85").Value = 1999.8334
$ws.Range("K85").Value = 1999.8334
$ws.Range("M85").Value = -751.8334

$ws.Range("H113").Value = 4911.05
$ws.Range("I113").Value = 4959.0527
$ws.Range("K113").Value = 4959.0527
$ws.Range("M113").Value = -2789.0527

$ws.Range("H126").Value = 2249.3333
$ws.Range("I126").Value = 1999.2
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 5997.6
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -3527.6
$ws.Range("N126").Value = -15440

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 133791.64
$ws.Range("I2").Value = 15454.728
$ws.Range("J2").Value = 350742.66
$ws.Range("K2").Value = 15454.728
$ws.Range("L2").Value = 350742.66
$ws.Range("M2").Value = -15342.728
$ws.Range("N2").Value = -350966.66

$ws.Range("H39").Value = 22024.5
$ws.Range("J39").Value = 22024.5
$ws.Range("L39").Value = 22024.5
$ws.Range("N39").Value = -22850.5

$ws.Range("H46").Value = 103579.1
$ws.Range("J46").Value = 103579.1
$ws.Range("L46").Value = 103579.1
$ws.Range("N46").Value = -104041.1

$ws.Range("H81").Value = 13476.8
$ws.Range("I81").Value = 5617.75
$ws.Range("J81").Value = 15805.407
$ws.Range("K81").Value = 11235.5
$ws.Range("L81").Value = 31610.814
$ws.Range("M81").Value = -10174.5
$ws.Range("N81").Value = -33732.814

$ws.Range("H84").Value = 13476.8
$ws.Range("I84").Value = 5617.75
$ws.Range("J84").Value = 15805.407
$ws.Range("K84").Value = 56177.5
$ws.Range("L84").Value = 158054.07
$ws.Range("M84").Value = -50873.5
$ws.Range("N84").Value = -168662.07

$ws.Range("H122").Value = 2708.386
$ws.Range("I122").Value = 2742.825
$ws.Range("J122").Value = 2627.353
$ws.Range("K122").Value = 8228.474999999999
$ws.Range("L122").Value = 7882.059
$ws.Range("M122").Value = -5778.474999999999
$ws.Range("N122").Value = -12782.059

$ws.Range("H134").Value = 103579.1
$ws.Range("J134").Value = 103579.1
$ws.Range("L134").Value = 310737.3
$ws.Range("N134").Value = -315807.3
